{"js": "// Minor refresh: add a new chapter on XYZ.\n// The document currently ends with an empty trailing paragraph; fill it in\n// with \"Chapter on XYZ \" and append a new paragraph with the body text\n// \"XYZ is the best thing\" right after it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The trailing (last) paragraph in the body is the empty one to fill in.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertText(\"Chapter on XYZ \", Word.InsertLocation.end);\n\n// Add the second new paragraph right after it.\nlastParagraph.insertParagraph(\"XYZ is the best thing\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Minor refresh: add a new chapter on XYZ.\n# The document currently ends with an empty trailing paragraph; fill it in\n# with \"Chapter on XYZ \" and append a new paragraph with the body text\n# \"XYZ is the best thing\" right after it.\n\n$d = $word.ActiveDocument\n\n$lastPara = $d.Paragraphs.Last\n$r = $lastPara.Range\n$r.InsertAfter(\"Chapter on XYZ \")\n\n# Create a new paragraph after the one we just filled in.\n$r.InsertParagraphAfter()\n\n# The newly created (empty) trailing paragraph is now the last paragraph;\n# fill it in with the second line of the new chapter.\n$newPara = $d.Paragraphs.Last\n$newPara.Range.InsertAfter(\"XYZ is the best thing\")\n"}
